$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10 through 26 (they get merged into rows 2-9)
$ws.Range("A10:A26").EntireRow.Delete() | Out-Null

$apos = [char]0x2019

$ws.Range("A2").Value = "('Beast', ['Token Creature — Beast', '3/3'])"
$ws.Range("A3").Value = "('Myr', ['Token Artifact Creature — Myr', '1/1'])"
$ws.Range("A4").Value = "('Rowan Kenrith Emblem', ['Emblem — Rowan', 'Whenever you activate an ability that isn" + $apos + "t a mana ability, copy it. You may choose new targets for the copy.'])"
$ws.Range("A5").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A6").Value = "('Warrior', ['Token Creature — Warrior', '1/1'])"
$ws.Range("A7").Value = "('Will Kenrith Emblem', ['Emblem — Will', 'Whenever you cast an instant or sorcery spell, copy it. You may choose new targets for the copy.'])"
$ws.Range("A8").Value = "('Zombie', ['Token Creature — Zombie', '2/2'])"
$ws.Range("A9").Value = "('Zombie Giant', ['Token Creature — Zombie Giant', '5/5'])"
